$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Re-style the "model input" columns (C:F) on rows 12-16 and 19-23 so they
#    all share the same look (font + orange-tinted fill) that a handful of
#    cells in rows 19-23 already used. Copy the existing formatting (from a
#    cell that already carries the target style) instead of re-describing it,
#    so the existing style/fill/font table entries are reused rather than
#    duplicated.
# ---------------------------------------------------------------------------
$ws.Range("C19").Copy() | Out-Null
$ws.Range("C12:F16").PasteSpecial(-4122) | Out-Null
$ws.Range("C19:F23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Refresh the underlying numbers (model re-run with new data).
# ---------------------------------------------------------------------------
$ws.Range("C12").Value = 9.15
$ws.Range("D12").Value = 7.09
$ws.Range("E12").Value = 13.23
$ws.Range("F12").Value = 29.24

$ws.Range("C13").Value = 1.65
$ws.Range("D13").Value = 3.14
$ws.Range("E13").Value = 5.3
$ws.Range("F13").Value = 7.47

$ws.Range("C14").Value = 686.17
$ws.Range("D14").Value = 484.28
$ws.Range("E14").Value = 833.86
$ws.Range("F14").Value = 2251.3000000000002

$ws.Range("C15").Value = 190.98
$ws.Range("D15").Value = 242.99
$ws.Range("E15").Value = 336.6
$ws.Range("F15").Value = 477.42

$ws.Range("C16").Value = 12.98
$ws.Range("D16").Value = 0.83
$ws.Range("E16").Value = 17.27
$ws.Range("F16").Value = 16.32

$ws.Range("C19").Value = 545.24
$ws.Range("D19").Value = 2345.58
$ws.Range("E19").Value = 3075.61
$ws.Range("F19").Value = 4818.71

$ws.Range("C20").Value = 267.33
$ws.Range("D20").Value = 1450.17
$ws.Range("E20").Value = 2930.08
$ws.Range("F20").Value = 3998.39

$ws.Range("C21").Value = 88.12
$ws.Range("D21").Value = 758.22
$ws.Range("E21").Value = 473.21
$ws.Range("F21").Value = 624.59

$ws.Range("C22").Value = 1557.66
$ws.Range("D22").Value = 6159.84
$ws.Range("E22").Value = 6899.08
$ws.Range("F22").Value = 6953.07

$ws.Range("C23").Value = 30.59
$ws.Range("D23").Value = 145.41
$ws.Range("E23").Value = 138.63999999999999
$ws.Range("F23").Value = 299.86

# ---------------------------------------------------------------------------
# 3) Move the active selection, as left by the author after the edit.
# ---------------------------------------------------------------------------
$ws.Range("B14").Select() | Out-Null
